$d = $word.ActiveDocument

# wdBrightGreen = 4 (maps to OOXML w:highlight w:val="green")
$green = 4

# Paragraph: "počet objektových vlastností: min. 7,"  (yellow -> green)
$d.Paragraphs(18).Range.Font.HighlightColorIndex = $green

# Paragraph: "ke každé objektové vlastnosti vytvořit vlastnost inverzní,"  (yellow -> green)
$d.Paragraphs(19).Range.Font.HighlightColorIndex = $green

# Paragraph: "vhodně použít logické konstruktory AND, OR nebo NOT (nepočítá se implicitní AND mezi jednotlivými logickými podmínkami uvnitř tříd),"  (none -> green)
$d.Paragraphs(23).Range.Font.HighlightColorIndex = $green

# Paragraph: "každá třída, která obsahuje jakékoliv omezení, musí být dle modelované skutečnosti řádně okomentována (v Protégé/Annotations),"  (none -> green)
$d.Paragraphs(26).Range.Font.HighlightColorIndex = $green
